# Userstories und Backlog update
# Johanna Terp und Constanze Richter

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Userstories")
$ws2 = $wb.Worksheets.Item("Konkretisierung")

# --- Sheet1 "Userstories": update cell texts (column positions unchanged) ---

# Header row
$ws1.Range("F4").Value = "Ziel"

# Herr Müller block
$ws1.Range("E5").Value = "braucht eine Karte der Messhalle"
$ws1.Range("F5").Value = "um Messpunkte visuell darstellen zu können"
$ws1.Range("G5").Value = "Ist eine Karte vorhanden, um Messpunkt visuell darstellen zu können"

$ws1.Range("E6").Value = "möchte einen neuen Messpunkt hinzufügen"
$ws1.Range("G6").Value = "An dieser Stelle ist ein Messpunkt hinzufügbar"

$ws1.Range("E7").Value = "möchte zu allen Messpunkten"
$ws1.Range("G7").Value = "Alle Messpunkte sind vorhanden um neue Messwerte hinzuzufügen"

# Frau Schneider block
$ws1.Range("E8").Value = "möchte Messpunkte archivieren"
$ws1.Range("F8").Value = "um nicht mehr benötigte Messpunkte abzulegen"

# row 9 unchanged

# Herr Leiser block
$ws1.Range("E10").Value = "möchte eine sortierte Ansicht der Grenzwertüberschreitungen"
$ws1.Range("F10").Value = "um Messpunkte bezogen auf ihr Datum nachzuvollziehen können"

$ws1.Range("E11").Value = "möchte eine sortierte Ansicht der Grenzwertüberschreitungen"
$ws1.Range("F11").Value = "um Messpunkte bezogen auf ihr Uhrzeit nachzuvollziehen können"

$ws1.Range("E12").Value = "möchte Messpunkte filtern können"
$ws1.Range("F12").Value = "um eine Auswahlmöglichkeit der Messpunkte zu bekommen"
$ws1.Range("G12").Value = "Messungen können nach Datum/ Uhrzeit gefiltert werden"

$ws1.Range("E13").Value = "möchte angelegte Daten exportieren"
$ws1.Range("F13").Value = "um sie für einen Geschäftsbericht zu nutzen"

$ws1.Range("E14").Value = "möchte angelegte Daten exportieren"
$ws1.Range("F14").Value = "um sie in externen Programmen visuelle darstellen zu können"

# Remove the leftover, empty formatted row 1 (no data, just a stray height)
$ws1.Rows.Item(1).EntireRow.AutoFit()

# --- Sheet views / active sheet ---
# Before: Konkretisierung (sheet2) was the active/selected tab.
# After: Userstories (sheet1) is active, with D17 selected;
# Konkretisierung keeps its own prior selection (L16) but is no longer the active tab.
$ws1.Activate()
$ws1.Range("D17").Select()
